$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "all": duplicate the last data row (36 -> 37) to add the new day's
# figures, pushing the footer note row down from 37 to 38. Row 36 is copied
# so the new row inherits its cell styles, then we overwrite the date.
# ---------------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("all")
$wsAll.Activate()
$wsAll.Rows("36").Copy()
$wsAll.Rows("37").Insert()
$wsAll.Range("A37").Value = 43965
$wsAll.Range("I37").Select()

# ---------------------------------------------------------------------------
# Sheet "kobe": duplicate the last data row (91 -> 92), pushing the footer
# note row down from 92 to 93, then fix up the values that differ from the
# prior day (date, and two counts back to 0).
# ---------------------------------------------------------------------------
$wsKobe = $wb.Worksheets.Item("kobe")
$wsKobe.Activate()
$wsKobe.Rows("91").Copy()
$wsKobe.Rows("92").Insert()
$wsKobe.Range("A92").Value = 43965
$wsKobe.Range("B92").Value = 0
$wsKobe.Range("D92").Value = 0
$wsKobe.Range("A92").Select()

# ---------------------------------------------------------------------------
# Sheet "other": duplicate the last data row (66 -> 67), pushing the footer
# note row down from 67 to 68, then fix up the date.
# ---------------------------------------------------------------------------
$wsOther = $wb.Worksheets.Item("other")
$wsOther.Activate()
$wsOther.Rows("66").Copy()
$wsOther.Rows("67").Insert()
$wsOther.Range("A67").Value = 43965
$wsOther.Range("A67").Select()

$wsAll.Activate()
